## Inserts a new data row at sheet row 152 (pushing existing rows 152-224
## down to 153-225) and populates it with a new "Espinaca" price record for
## "Vega Modelo de Temuco" / La Araucanía.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 152, shifting everything below down.
$ws.Rows.Item(152).Insert()

# Fill in the new row 152 with the record's data.
$ws.Cells.Item(152, 1).Value2  = 10
$ws.Cells.Item(152, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(152, 3).Value2  = "La Araucanía"
$ws.Cells.Item(152, 4).Value2  = 44917
$ws.Cells.Item(152, 5).Value2  = 9
$ws.Cells.Item(152, 6).Value2  = 100112012
$ws.Cells.Item(152, 7).Value2  = "Espinaca"
$ws.Cells.Item(152, 8).Value2  = "Sin especificar"
$ws.Cells.Item(152, 9).Value2  = "Primera"
$ws.Cells.Item(152, 10).Value2 = 55
$ws.Cells.Item(152, 11).Value2 = 10000
$ws.Cells.Item(152, 12).Value2 = 10000
$ws.Cells.Item(152, 13).Value2 = 10000
$ws.Cells.Item(152, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(152, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(152, 16).Value2 = 3333
$ws.Cells.Item(152, 17).Value2 = 3
$ws.Cells.Item(152, 18).Value2 = "Hortaliza"
